# Add the "unit"-style columns (D:G) for each data row, and update the
# active selection to match the new range the author was working in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "U"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1

# Row 3
$ws.Range("D3").Value = "U"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 2

# Row 4
$ws.Range("D4").Value = "U"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2

# Row 5
$ws.Range("D5").Value = "U"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 3

# Match the author's final selection (D2 active cell, D2:G5 selected).
$ws.Range("D2:G5").Select()
